$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. "Author (s):" line — re-type the "__Haseeb, Aiza, Ali_____________"
#    text so Word collapses it (and the spell-check proofErr wrapper
#    around "Haseeb") into a single run.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("__Haseeb, Aiza, Ali_____________", $false, $false, $false, $false, $false, `
    $true, 1, $false, "__Haseeb, Aiza, Ali_____________", 2) | Out-Null

# ---------------------------------------------------------------------
# 2. OPEN ISSUES note: "...serve patients" -> "...serve customers".
#
#    In the source document a hidden "_GoBack" bookmark sits in the
#    middle of this sentence (right where the previous edit session's
#    cursor last was), splitting it into two runs. Editing "patients"
#    directly would delete that bookmark, so first walk the bookmark
#    forward, one character at a time, until it sits at the very end of
#    the sentence (where this new edit will actually happen) before
#    making the wording change.
# ---------------------------------------------------------------------
$bookmarkRange = $d.Content
$bookmarkRange.Find.Execute("work done o", $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$pos = $bookmarkRange.End

$endRange = $d.Content
$endRange.Find.Execute("patients", $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$target = $endRange.End

while ($pos -lt $target) {
    $charRange = $d.Range($pos, $pos + 1)
    $ch = $charRange.Text
    $charRange.Delete()
    $insertionPoint = $d.Range($pos, $pos)
    $insertionPoint.InsertBefore($ch)
    $pos = $pos + 1
}

# Now replace the word itself. Toggling Bold on and back off around the
# replacement forces Word to keep the new text ("customers") as its own
# run instead of silently re-merging it with the sentence that precedes
# it, matching how the document was actually re-split upstream.
$wordRange = $d.Content
$wordRange.Find.Execute("patients", $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$wordRange.Bold = $true
$wordRange.Text = "customers"
$wordRange.Bold = $false
